$wb = $excel.ActiveWorkbook

# Rename sheets to reflect the unified DataNode / DataTable / Entity concept
$wsDataNode = $wb.Worksheets.Item("Property1")
$wsDataNode.Name = "DataNode"

$wsDataTable = $wb.Worksheets.Item("Record")
$wsDataTable.Name = "DataTable"

# Make the DataTable sheet the active/selected sheet (tab) on reopen
$wsDataTable.Activate()
